$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:D51 is treated as text so numeric-looking values (e.g. "1.00", "39.90")
# keep their exact original formatting instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '58.536.03'
$ws.Range("E2").Value = '  +4.15%  '
$ws.Range("D3").Value = '3.308.36'
$ws.Range("E3").Value = '  +2.65%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '401.25'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").Value = '110.65'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '0.585'
$ws.Range("E7").Value = '  +5.85%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.631'
$ws.Range("E9").Value = '  +2.24%  '
$ws.Range("D10").Value = '39.90'
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("D11").Value = '0.0971'
$ws.Range("E11").Value = '  +5.87%  '
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").Value = '3.819.87'
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("D14").Value = '8.41'
$ws.Range("E14").Value = '  +4.24%  '
$ws.Range("D15").Value = '19.17'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '3.299.52'
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '11.06'
$ws.Range("E18").Value = '  +2.37%  '
$ws.Range("D19").Value = '58.199.72'
$ws.Range("E19").Value = '  +3.96%  '
$ws.Range("D20").Value = '3.35'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = '0.0000109'
$ws.Range("E21").Value = '  +5.75%  '
$ws.Range("D22").Value = '13.06'
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("D23").Value = '300.09'
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").Value = '74.96'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").Value = '3.20'
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").Value = '28.40'
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("D27").Value = '8.00'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").Value = '4.43'
$ws.Range("E28").Value = '  +1.35%  '
$ws.Range("D29").Value = '7.35'
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '0.113'
$ws.Range("E31").Value = '  +2.65%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("D33").Value = '11.35'
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D34").Value = '41.42'
$ws.Range("E34").Value = '  +14.24%  '
$ws.Range("D35").Value = '0.0507'
$ws.Range("E35").Value = '  +2.48%  '
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("D37").Value = '51.80'
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("D38").Value = '3.29'
$ws.Range("E38").Value = '  +4.86%  '
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").Value = '0.997'
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").Value = '138.32'
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  +2.45%  '
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '3.95'
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").Value = '16.97'
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '0.282'
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +8.25%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '22.58'
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("D49").Value = '2.174.11'
$ws.Range("E49").Value = '  +2.66%  '
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("D51").Value = '1.94'
$ws.Range("E51").Value = '  -11.31%  '

# Restore the original (default) cell style on the Price column now that the
# text values have been written, so formatting matches the source workbook.
$ws.Range("D2:D51").Style = "Normal"

